$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @(
    "Mon Apr 08 17:33:04 EDT 2024",
    "Mon Apr 08 17:33:13 EDT 2024",
    "Mon Apr 08 17:33:23 EDT 2024",
    "Mon Apr 08 17:33:33 EDT 2024",
    "Mon Apr 08 17:33:43 EDT 2024",
    "Mon Apr 08 17:33:53 EDT 2024",
    "Mon Apr 08 17:34:02 EDT 2024",
    "Mon Apr 08 17:34:12 EDT 2024",
    "Mon Apr 08 17:34:22 EDT 2024",
    "Mon Apr 08 17:34:31 EDT 2024",
    "Mon Apr 08 17:34:41 EDT 2024",
    "Mon Apr 08 17:34:51 EDT 2024",
    "Mon Apr 08 17:35:00 EDT 2024",
    "Mon Apr 08 17:35:10 EDT 2024",
    "Mon Apr 08 17:35:20 EDT 2024",
    "Mon Apr 08 17:35:30 EDT 2024",
    "Mon Apr 08 17:35:39 EDT 2024",
    "Mon Apr 08 17:35:49 EDT 2024",
    "Mon Apr 08 17:35:59 EDT 2024",
    "Mon Apr 08 17:36:08 EDT 2024",
    "Mon Apr 08 17:36:18 EDT 2024",
    "Mon Apr 08 17:36:28 EDT 2024",
    "Mon Apr 08 17:36:37 EDT 2024",
    "Mon Apr 08 17:36:47 EDT 2024",
    "Mon Apr 08 17:36:57 EDT 2024",
    "Mon Apr 08 17:37:06 EDT 2024",
    "Mon Apr 08 17:37:16 EDT 2024",
    "Mon Apr 08 17:37:26 EDT 2024",
    "Mon Apr 08 17:37:35 EDT 2024"
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $dates[$i]
}
